$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (shifts old F "District" to G)
$ws.Columns.Item(6).Insert()

# Fill in new column F with Address values
$ws.Range("F2").Value = "Address"
$ws.Range("F3").Value = "Govt. Science College"
$ws.Range("F4").Value = "G H S kampalapuraPiriyapattana Taluk"
$ws.Range("F5").Value = "T G T G H P S DegganahalliK R Nagar"
$ws.Range("F6").Value = "Sri Vani Vilasa Ursu Girls High SchoolNazarbad"
$ws.Range("F7").Value = "G H P S Kudanahally"
$ws.Range("F8").Value = "Govt. High SchoolBettadatungaPeriyapatna"
$ws.Range("F9").Value = "Sri Royal High SchoolD Saluhundi"
$ws.Range("F10").Value = "G H S KedagaK R Nagar"
$ws.Range("F11").Value = "Dayananda High School Bherya K R Nagar"
$ws.Range("F12").Value = "J S S High School ManuganahalliH D Kote Taluk"
$ws.Range("F13").Value = "S C V D S High School Bettadapur"
$ws.Range("F14").Value = "G H S MaggeH D Kote"
$ws.Range("F15").Value = "Govt. High SchoolD Salundi"
$ws.Range("F16").Value = "G H S KolathurT Narasipura"
$ws.Range("F17").Value = "Sacred Heart High SchoolNaganahallyH D Kote"
$ws.Range("F18").Value = "Vidya Vardhaka High SchoolMandi Mohalla"
$ws.Range("F19").Value = "TGT G H P SBharathawadiHunsur"
$ws.Range("F20").Value = "Sharadavilas High School"
$ws.Range("F21").Value = "G H S Devalapura"
$ws.Range("F22").Value = "Govt High School Antharasanthe H D Kote Taluk"
$ws.Range("F23").Value = "G H S Siddaramana Hundi"
$ws.Range("F24").Value = "Marimallappa High School"
$ws.Range("F25").Value = "Sri Vidyadayini girls high school"
$ws.Range("F26").Value = "S N S High School GavadagereHunsur"
$ws.Range("F27").Value = "G H S ChikkahunsurHunsurTaluk"
$ws.Range("F28").Value = "Govt. Junior College (High School Section) Yelwala"
$ws.Range("F29").Value = "Sri Valmiki High School ChallahalliHunsur"
$ws.Range("F30").Value = "Govt. Junior College for BoysHunsur"
$ws.Range("F31").Value = "G P U C Hediyala Nanjangud"
$ws.Range("F32").Value = "M R P High SchoolAralimarada Koppal Periyapatna"
$ws.Range("F33").Value = "G G H ST Narasipura"
$ws.Range("F34").Value = "Sri MaruthiHigh School ThattekereHunsur"
$ws.Range("F35").Value = "G H P S Lakshmipura K R Nagar"
$ws.Range("F36").Value = "G H S Harohally Mellahally"
$ws.Range("F37").Value = "G J C KitturPeriyapatna"
$ws.Range("F38").Value = "Sri K Puttaswamy High School Kumbarakoppal"
$ws.Range("F39").Value = "D Kote Taluk"
$ws.Range("F40").Value = "Shree Gurumallewara High School karyaNanjfud"
$ws.Range("F41").Value = "G J C RavandurPeriyapatna"
$ws.Range("F42").Value = "Govt. High School TadimalangiT Narasipura (tq)"
$ws.Range("F43").Value = "Avila ConventHigh School"
$ws.Range("F44").Value = "G H S HedathaleNanjangud"
$ws.Range("F45").Value = "Govt. High SchoolVataluT Narasipura"
$ws.Range("F46").Value = "G H S MusuvinakoppaluT N Pura"
$ws.Range("F47").Value = "Bhagini Seva Samaja High SchoolK M Puram"
$ws.Range("F48").Value = "Dayananda High School DevithandreK R Nagar"
$ws.Range("F49").Value = "G H S HeggadahallyNanjanagud"
$ws.Range("F50").Value = "Sri Kuvempu High School Kuvempu Nagar"
$ws.Range("F51").Value = "G H S HeggurT N Pura"
$ws.Range("F52").Value = "G J C Periyapettna TqHalaganahally"
$ws.Range("F53").Value = "ManasagangothriHigh SchoolManasagangothri"
